$d = $word.ActiveDocument

# 1. Title heading (Heading1) and later bold repeat near the end both share
#    the exact same original text, so replace every occurrence in the body.
$d.Content.Find.Execute(
    "Play Day of the Dead for Free " + [char]0x2013 + " Unique Gameplay Mechanics",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Day of the Dead Slot for Free", 2)

# 2. "What we like" bullet list items
$d.Content.Find.Execute(
    "Detailed graphics and sound effects",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Beneficial scatter, wild, and bonus symbols", 2)

$d.Content.Find.Execute(
    "Mexican-inspired theme is well-executed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impressive graphics and sound design", 2)

$d.Content.Find.Execute(
    "720 possible payline combinations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Mexican-inspired theme and symbols", 2)

# 3. "What we don't like" bullet list items
$d.Content.Find.Execute(
    "Possibility of fewer major paylines",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Possibility of hitting major paylines less frequently", 2)

$d.Content.Find.Execute(
    "Dia de Los Muertos theme may not appeal to all players",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limited availability of slots with the same theme", 2)

# 4. Meta title/description paragraphs near the end of the document
$d.Content.Find.Execute(
    "Try Day of the Dead, a great online slot game based on the Mexican celebration. Enjoy unique gameplay mechanics, great graphics, and sound. Play for free today.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Experience the unique gameplay and Mexican-inspired theme of Day of the Dead slot for free.", 2)
